$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Completed")

# Row 10: The First 20 Minutes
$ws.Range("A10").Value = "The First 20 Minutes"
$ws.Range("B10").Value = "Gretchen Reynolds"

$ws.Range("C9").Copy()
$ws.Range("C10:D10").PasteSpecial(-4122)
$ws.Range("C10").Value = 43843
$ws.Range("D10").Value = 43846

$ws.Range("E10").Value = "exercise;science;biology;health;fitness"
$ws.Range("F10").Value = "Audio"
$ws.Range("G10").Value = "9 Hrs 7 Mins"

# Row 11: Mindset
$ws.Range("A11").Value = "Mindset"
$ws.Range("B11").Value = "Carol Dweck"

$ws.Range("C9").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$ws.Range("C11").Value = 43842

$ws.Range("E11").Value = "psychology;self-improvement;growth mindset;fixed mindset"
$ws.Range("F11").Value = "Hard Copy"
$ws.Range("G11").Value = "320 Pages"

$excel.CutCopyMode = 0
$ws.Range("A12").Select()
